# "Generate Report for handoff"
#
# The two localization records (12cca4fe... and d993f17d...) swap their
# display order (d993f17d now listed first / row 2, 12cca4fe now listed
# second / row 3) on every sheet, and the 12cca4fe record is marked ready
# for a new handoff (status + handoff datetime refreshed) since it was
# just regenerated/handed off again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
$ws.Range("A3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$hlMap = @{
    '$A$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
    '$A$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
}
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($hlMap.ContainsKey($addr)) {
        $h.TextToDisplay = $hlMap[$addr]
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
$ws.Range("C2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.zh-cn.xlf"
$ws.Range("E2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
$ws.Range("F2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.zh-cn.xlf"

$ws.Range("A3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.zh-cn.xlf"
$ws.Range("D3").Value = "2016-02-16 14:10:59"
$ws.Range("E3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
$ws.Range("F3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.zh-cn.xlf"

$hlMap = @{
    '$A$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
    '$C$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.zh-cn.xlf"
    '$E$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
    '$F$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.zh-cn.xlf"
    '$A$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
    '$C$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.zh-cn.xlf"
    '$E$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
    '$F$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.zh-cn.xlf"
}
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($hlMap.ContainsKey($addr)) {
        $h.TextToDisplay = $hlMap[$addr]
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
$ws.Range("C2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.de-de.xlf"
$ws.Range("E2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
$ws.Range("F2").Value = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.de-de.xlf"

$ws.Range("A3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.de-de.xlf"
$ws.Range("D3").Value = "2016-02-16 14:11:12"
$ws.Range("E3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
$ws.Range("F3").Value = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.de-de.xlf"

$hlMap = @{
    '$A$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
    '$C$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.de-de.xlf"
    '$E$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.md"
    '$F$2' = "d993f17d-ff44-4234-9dfc-834532cfa995.7053f1e9c8f10e31af3a529038d52d63194e6683.de-de.xlf"
    '$A$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
    '$C$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.de-de.xlf"
    '$E$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.md"
    '$F$3' = "12cca4fe-e9f2-4889-aa2f-686374696b5a.38ddf27578073a1fdae0406cd0ef4058b2f119d3.de-de.xlf"
}
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($hlMap.ContainsKey($addr)) {
        $h.TextToDisplay = $hlMap[$addr]
    }
}
